$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioA")
$ws.Range("I1").Value = 200
$ws.Activate()
$ws.Range("I2").Select()
$excel.CalculateFullRebuild()
